# Registration files/Service_list.xlsx - add new service log rows (error handling logs, WIP)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: VIN (column D) becomes a real number instead of text ---
$ws.Cells.Item(4, 4).Value = 12345678912345680

# --- Row 5 (new) ---
$ws.Cells.Item(5, 1).Value = 8138074349
$ws.Cells.Item(5, 2).Value = "Yo yo"
$ws.Cells.Item(5, 3).Value = 992907510905
$ws.Cells.Item(5, 4).Value = 12345678912345680
$ws.Cells.Item(5, 5).Value = "Camry 10"
$ws.Cells.Item(5, 6).Value = "Service 3"
$ws.Cells.Item(5, 7).Value = "13/08/2025"
$ws.Cells.Item(5, 8).Value = "2025-08-11 09:30:28"
$ws.Cells.Item(5, 9).Value = "09:30"
$ws.Cells.Item(5, 10).Value = "'12345"

# --- Row 6 (new) ---
$ws.Cells.Item(6, 1).Value = 8138074349
$ws.Cells.Item(6, 2).Value = "Yo yo"
$ws.Cells.Item(6, 3).Value = 992907510905
$ws.Cells.Item(6, 4).Value = 12354678912345680
$ws.Cells.Item(6, 5).Value = "'123"
$ws.Cells.Item(6, 6).Value = "Service 1"
$ws.Cells.Item(6, 7).Value = "16/08/2025"
$ws.Cells.Item(6, 8).Value = "2025-08-11 09:35:39"
$ws.Cells.Item(6, 9).Value = "09:30"
$ws.Cells.Item(6, 10).Value = "-"

# --- Row 7 (new) ---
$ws.Cells.Item(7, 1).Value = 8138074349
$ws.Cells.Item(7, 2).Value = "Yo yo"
$ws.Cells.Item(7, 3).Value = 992907510905
$ws.Cells.Item(7, 4).Value = 12345678912345680
$ws.Cells.Item(7, 5).Value = "Toyota Camry 6"
$ws.Cells.Item(7, 6).Value = "Service 2"
$ws.Cells.Item(7, 7).Value = "13/08/2025"
$ws.Cells.Item(7, 8).Value = "2025-08-11 10:47:28"
$ws.Cells.Item(7, 9).Value = "15:00"
$ws.Cells.Item(7, 10).Value = "-"

# --- Row 8 (new) - keeps VIN as text (quote-prefixed) like the original row 4 did ---
$ws.Cells.Item(8, 1).Value = 8138074349
$ws.Cells.Item(8, 2).Value = "Yo yo"
$ws.Cells.Item(8, 3).Value = 992907510905
$ws.Cells.Item(8, 4).Value = "'12345678912345678"
$ws.Cells.Item(8, 5).Value = "Toyota Camry 6"
$ws.Cells.Item(8, 6).Value = "Service 3"
$ws.Cells.Item(8, 7).Value = "30/08/2025"
$ws.Cells.Item(8, 8).Value = "2025-08-11 12:03:40"
$ws.Cells.Item(8, 9).Value = "15:00"
$ws.Cells.Item(8, 10).Value = "Nope"
